# Updated data dictionary acc table
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ACC")

# Make sure ACC is the active/selected sheet (it already is tabSelected in the
# source file, but Activate() keeps things consistent).
$ws.Activate()

# --- Row 3 : Address field -------------------------------------------------
# B3 already carries a style (s="1"); just give it the "Address" text.
$ws.Range("B3").Value = "Address"

# --- Row 4 : Phone field (new row) -----------------------------------------
# New shared strings must be introduced in the same order they appear in the
# target workbook's sharedStrings table, so set these in this exact order:
#   186 Phone
#   187 Refers to the injury codes
#   188 The offices phone number
#   189 The offices address
$ws.Range("B4").Value = "Phone"

# --- Row 5 : Acc Code ID field ----------------------------------------------
$ws.Range("C5").Value = "Refers to the injury codes"

# --- finish Row 4 ------------------------------------------------------------
$ws.Range("C4").Value = "The offices phone number"

# --- finish Row 3 ------------------------------------------------------------
$ws.Range("C3").Value = "The offices address"

# --- remaining (re-used) shared strings + numeric values --------------------
$ws.Range("D3").Value = "Varchar"
$ws.Range("E3").Value = 30

$ws.Range("D4").Value = "Integer"
$ws.Range("E4").Value = 10

$ws.Range("B5").Value = "Acc Code ID"

# --- selection change reflected in the sheet view ---------------------------
$ws.Range("C6").Select()
